$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.73148533333334
$ws.Range("H2").Value = 71.194456
$ws.Range("I2").Value = 0.8653076146801144
$ws.Range("J2").Value = 0.8653076146801145
$ws.Range("M2").Value = 8.813278666666667
$ws.Range("N2").Value = 26.439836
$ws.Range("O2").Value = 0.3770976991891536
$ws.Range("P2").Value = 0.3770976991891536
$ws.Range("Q2").Value = 209.1521934165796
$ws.Range("R2").Value = 1882.369740749216
$ws.Range("S2").Value = 0.3263055105867259
$ws.Range("T2").Value = 0.3263055105867259
$ws.Range("G3").Value = 23.73148533333334
$ws.Range("H3").Value = 71.194456
$ws.Range("I3").Value = 0.8653076146801144
$ws.Range("J3").Value = 0.8653076146801145
$ws.Range("O3").Value = 0.5522024902836482
$ws.Range("P3").Value = 0.5522024902836482
$ws.Range("Q3").Value = 306.2717229547192
$ws.Range("R3").Value = 2756.445506592472
$ws.Range("S3").Value = 0.4778250196877627
$ws.Range("T3").Value = 0.4778250196877628
$ws.Range("G4").Value = 23.73148533333334
$ws.Range("H4").Value = 71.194456
$ws.Range("I4").Value = 0.8653076146801144
$ws.Range("J4").Value = 0.8653076146801145
$ws.Range("M4").Value = 1.649921333333333
$ws.Range("N4").Value = 4.949764
$ws.Range("O4").Value = 0.07059592260441032
$ws.Range("P4").Value = 0.07059592260441033
$ws.Range("Q4").Value = 39.15508392315378
$ws.Range("R4").Value = 352.395755308384
$ws.Range("S4").Value = 0.06108718939496426
$ws.Range("T4").Value = 0.06108718939496428
$ws.Range("G5").Value = 23.73148533333334
$ws.Range("H5").Value = 71.194456
$ws.Range("I5").Value = 0.8653076146801144
$ws.Range("J5").Value = 0.8653076146801145
$ws.Range("M5").Value = 0.002428
$ws.Range("N5").Value = 0.007284
$ws.Range("O5").Value = 0.0001038879227879399
$ws.Range("P5").Value = 0.0001038879227879399
$ws.Range("Q5").Value = 0.05762004638933334
$ws.Range("R5").Value = 0.518580417504
$ws.Range("S5").Value = 0.00008989501066170422
$ws.Range("T5").Value = 0.00008989501066170423
$ws.Range("I6").Value = 0.09010639372350319
$ws.Range("J6").Value = 0.09010639372350321
$ws.Range("M6").Value = 8.813278666666667
$ws.Range("N6").Value = 26.439836
$ws.Range("O6").Value = 0.3770976991891536
$ws.Range("P6").Value = 0.3770976991891536
$ws.Range("Q6").Value = 21.77948000041067
$ws.Range("R6").Value = 196.015320003696
$ws.Range("S6").Value = 0.03397891375536505
$ws.Range("T6").Value = 0.03397891375536505
$ws.Range("I7").Value = 0.09010639372350319
$ws.Range("J7").Value = 0.09010639372350321
$ws.Range("O7").Value = 0.5522024902836482
$ws.Range("P7").Value = 0.5522024902836482
$ws.Range("S7").Value = 0.04975697500459735
$ws.Range("T7").Value = 0.04975697500459736
$ws.Range("I8").Value = 0.09010639372350319
$ws.Range("J8").Value = 0.09010639372350321
$ws.Range("M8").Value = 1.649921333333333
$ws.Range("N8").Value = 4.949764
$ws.Range("O8").Value = 0.07059592260441032
$ws.Range("P8").Value = 0.07059592260441033
$ws.Range("Q8").Value = 4.077305397989333
$ws.Range("R8").Value = 36.695748581904
$ws.Range("S8").Value = 0.006361143997466955
$ws.Range("T8").Value = 0.006361143997466958
$ws.Range("I9").Value = 0.09010639372350319
$ws.Range("J9").Value = 0.09010639372350321
$ws.Range("M9").Value = 0.002428
$ws.Range("N9").Value = 0.007284
$ws.Range("O9").Value = 0.0001038879227879399
$ws.Range("P9").Value = 0.0001038879227879399
$ws.Range("Q9").Value = 0.006000102736
$ws.Range("R9").Value = 0.054000924624
$ws.Range("S9").Value = 0.000009360966073847017
$ws.Range("T9").Value = 0.000009360966073847018
$ws.Range("G10").Value = 1.148663
$ws.Range("H10").Value = 3.445989
$ws.Range("I10").Value = 0.04188304383987305
$ws.Range("J10").Value = 0.04188304383987305
$ws.Range("M10").Value = 8.813278666666667
$ws.Range("N10").Value = 26.439836
$ws.Range("O10").Value = 0.3770976991891536
$ws.Range("P10").Value = 0.3770976991891536
$ws.Range("Q10").Value = 10.12348711308933
$ws.Range("R10").Value = 91.111384017804
$ws.Range("S10").Value = 0.01579399946705458
$ws.Range("T10").Value = 0.01579399946705458
$ws.Range("G11").Value = 1.148663
$ws.Range("H11").Value = 3.445989
$ws.Range("I11").Value = 0.04188304383987305
$ws.Range("J11").Value = 0.04188304383987305
$ws.Range("O11").Value = 0.5522024902836482
$ws.Range("P11").Value = 0.5522024902836482
$ws.Range("Q11").Value = 14.82431424594367
$ws.Range("R11").Value = 133.418828213493
$ws.Range("S11").Value = 0.02312792110903711
$ws.Range("T11").Value = 0.02312792110903711
$ws.Range("G12").Value = 1.148663
$ws.Range("H12").Value = 3.445989
$ws.Range("I12").Value = 0.04188304383987305
$ws.Range("J12").Value = 0.04188304383987305
$ws.Range("M12").Value = 1.649921333333333
$ws.Range("N12").Value = 4.949764
$ws.Range("O12").Value = 0.07059592260441032
$ws.Range("P12").Value = 0.07059592260441033
$ws.Range("Q12").Value = 1.895203588510667
$ws.Range("R12").Value = 17.056832296596
$ws.Range("S12").Value = 0.002956772121356802
$ws.Range("T12").Value = 0.002956772121356803
$ws.Range("G13").Value = 1.148663
$ws.Range("H13").Value = 3.445989
$ws.Range("I13").Value = 0.04188304383987305
$ws.Range("J13").Value = 0.04188304383987305
$ws.Range("M13").Value = 0.002428
$ws.Range("N13").Value = 0.007284
$ws.Range("O13").Value = 0.0001038879227879399
$ws.Range("P13").Value = 0.0001038879227879399
$ws.Range("Q13").Value = 0.002788953764
$ws.Range("R13").Value = 0.025100583876
$ws.Range("S13").Value = 0.000004351142424560635
$ws.Range("T13").Value = 0.000004351142424560635
$ws.Range("G14").Value = 0.07412966666666666
$ws.Range("H14").Value = 0.222389
$ws.Range("I14").Value = 0.002702947756509242
$ws.Range("J14").Value = 0.002702947756509243
$ws.Range("M14").Value = 8.813278666666667
$ws.Range("N14").Value = 26.439836
$ws.Range("O14").Value = 0.3770976991891536
$ws.Range("P14").Value = 0.3770976991891536
$ws.Range("Q14").Value = 0.6533254098004444
$ws.Range("R14").Value = 5.879928688204
$ws.Range("S14").Value = 0.00101927538000812
$ws.Range("T14").Value = 0.00101927538000812
$ws.Range("G15").Value = 0.07412966666666666
$ws.Range("H15").Value = 0.222389
$ws.Range("I15").Value = 0.002702947756509242
$ws.Range("J15").Value = 0.002702947756509243
$ws.Range("O15").Value = 0.5522024902836482
$ws.Range("P15").Value = 0.5522024902836482
$ws.Range("Q15").Value = 0.9566961533658889
$ws.Range("R15").Value = 8.610265380293001
$ws.Range("S15").Value = 0.001492574482251003
$ws.Range("T15").Value = 0.001492574482251004
$ws.Range("G16").Value = 0.07412966666666666
$ws.Range("H16").Value = 0.222389
$ws.Range("I16").Value = 0.002702947756509242
$ws.Range("J16").Value = 0.002702947756509243
$ws.Range("M16").Value = 1.649921333333333
$ws.Range("N16").Value = 4.949764
$ws.Range("O16").Value = 0.07059592260441032
$ws.Range("P16").Value = 0.07059592260441033
$ws.Range("Q16").Value = 0.1223081184662222
$ws.Range("R16").Value = 1.100773066196
$ws.Range("S16").Value = 0.000190817090622291
$ws.Range("T16").Value = 0.000190817090622291
$ws.Range("G17").Value = 0.07412966666666666
$ws.Range("H17").Value = 0.222389
$ws.Range("I17").Value = 0.002702947756509242
$ws.Range("J17").Value = 0.002702947756509243
$ws.Range("M17").Value = 0.002428
$ws.Range("N17").Value = 0.007284
$ws.Range("O17").Value = 0.0001038879227879399
$ws.Range("P17").Value = 0.0001038879227879399
$ws.Range("Q17").Value = 0.0001799868306666667
$ws.Range("R17").Value = 0.001619881476
$ws.Range("S17").Value = 0.0000002808036278280677
$ws.Range("T17").Value = 0.0000002808036278280677

Write-Host "Updated 174 cells"